# Auto-generated edit script: applies per-cell text updates to Sheet1
# matching the authoritative diff (cryptos.xlsx price/volume refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.991.92"
$ws.Range("E2").Value = "  -1.41%  "
# Row 3
$ws.Range("D3").Value = "3.178.34"
$ws.Range("E3").Value = "  +1.33%  "
# Row 4
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").Value = "'589.47"
$ws.Range("E5").Value = "  -1.97%  "
# Row 6
$ws.Range("D6").Value = "'138.57"
# Row 8
$ws.Range("D8").Value = "3.175.31"
$ws.Range("E8").Value = "  +1.44%  "
# Row 9
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  -1.08%  "
# Row 10
$ws.Range("E10").Value = "  -1.78%  "
# Row 11
$ws.Range("D11").Value = "'5.37"
$ws.Range("E11").Value = "  -0.10%  "
# Row 12
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  -1.31%  "
# Row 13
$ws.Range("E13").Value = "  -3.61%  "
# Row 14
$ws.Range("D14").Value = "'34.02"
$ws.Range("E14").Value = "  -2.90%  "
# Row 15
$ws.Range("D15").Value = "3.696.93"
$ws.Range("E15").Value = "  +1.30%  "
# Row 16
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "  +0.71%  "
# Row 17
$ws.Range("D17").Value = "3.171.83"
$ws.Range("E17").Value = "  +1.51%  "
# Row 18
$ws.Range("D18").Value = "62.998.62"
$ws.Range("E18").Value = "  -1.47%  "
# Row 19
$ws.Range("D19").Value = "'6.68"
$ws.Range("E19").Value = "  -2.00%  "
# Row 20
$ws.Range("D20").Value = "'473.20"
$ws.Range("E20").Value = "  -2.44%  "
# Row 21
$ws.Range("D21").Value = "'13.95"
$ws.Range("E21").Value = "  -4.82%  "
# Row 22
$ws.Range("D22").Value = "'0.703"
$ws.Range("E22").Value = "  -0.58%  "
# Row 23
$ws.Range("D23").Value = "'7.74"
$ws.Range("E23").Value = "  +1.44%  "
# Row 24
$ws.Range("D24").Value = "'83.64"
$ws.Range("E24").Value = "  -3.77%  "
# Row 25
$ws.Range("D25").Value = "'13.01"
$ws.Range("E25").Value = "  -3.23%  "
# Row 27
$ws.Range("D27").Value = "'2.71"
$ws.Range("E27").Value = "  -1.57%  "
# Row 28
$ws.Range("D28").Value = "'7.07"
$ws.Range("E28").Value = "  +1.00%  "
# Row 29
$ws.Range("D29").Value = "'7.95"
$ws.Range("E29").Value = "  -3.81%  "
# Row 30
$ws.Range("E30").Value = "  +1.06%  "
# Row 31
$ws.Range("E31").Value = "  -0.05%  "
# Row 32
$ws.Range("D32").Value = "'26.88"
$ws.Range("E32").Value = "  -0.88%  "
# Row 33
$ws.Range("E33").Value = "  -3.11%  "
# Row 34
$ws.Range("D34").Value = "'2.52"
$ws.Range("E34").Value = "  -4.67%  "
# Row 35
$ws.Range("D35").Value = "'1.08"
$ws.Range("E35").Value = "  -2.47%  "
# Row 36
$ws.Range("D36").Value = "'52.51"
$ws.Range("E36").Value = "  -0.05%  "
# Row 37
$ws.Range("D37").Value = "'5.78"
$ws.Range("E37").Value = "  -3.60%  "
# Row 38
$ws.Range("D38").Value = "0.0₃0707"
$ws.Range("E38").Value = "  -4.99%  "
# Row 39
$ws.Range("D39").Value = "'0.0388"
$ws.Range("E39").Value = "  -1.89%  "
# Row 40
$ws.Range("D40").Value = "'419.04"
$ws.Range("E40").Value = "  -4.42%  "
# Row 41
$ws.Range("D41").Value = "2.961.15"
$ws.Range("E41").Value = "  +2.94%  "
# Row 42
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.32"
$ws.Range("E42").Value = "  +0.44%  "
# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.73"
$ws.Range("E43").Value = "  -7.75%  "
# Row 44
$ws.Range("D44").Value = "'0.111"
$ws.Range("E44").Value = "  -7.86%  "
# Row 45
$ws.Range("E45").Value = "  +1.02%  "
# Row 46
$ws.Range("E46").Value = "  +0.06%  "
# Row 47
$ws.Range("D47").Value = "'2.13"
$ws.Range("E47").Value = "  -3.42%  "
# Row 48
$ws.Range("D48").Value = "'25.54"
$ws.Range("E48").Value = "  -1.44%  "
# Row 49
$ws.Range("E49").Value = "  -0.02%  "
# Row 50
$ws.Range("E50").Value = "  -5.75%  "
# Row 51
$ws.Range("D51").Value = "'119.52"
$ws.Range("E51").Value = "  -1.51%  "
